{"js": "// The cover-page template has two date FORMTEXT fields (Ausbildungsbeginn /\n// \"training start\" and Ausbildungsende / \"training end\") whose cached field\n// result still holds the placeholder text the template ships with\n// (\"00.00.0001\" / \"00.00.0002\", split across two <w:r> runs each). The fix\n// replaces that leftover placeholder with real sample dates.\n//\n// We locate each placeholder by searching the document body for its exact\n// text (this matches across the run boundary and returns one Range per\n// hit), then replace the matched range's text in place so the field keeps\n// its FORMTEXT structure (begin/separate/end field-chars + bookmark) and\n// ends up with a single run holding the new date.\n\nconst replacements = [\n  { find: \"00.00.0001\", replaceWith: \"11.11.1970\" },\n  { find: \"00.00.0002\", replaceWith: \"22.22.1970\" }\n];\n\nfor (const { find, replaceWith } of replacements) {\n  const hits = context.document.body.search(find, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(replaceWith, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The cover-page template has two date FORMTEXT fields (Ausbildungsbeginn /\n# \"training start\" and Ausbildungsende / \"training end\") whose cached field\n# result still holds the placeholder text the template ships with\n# (\"00.00.0001\" / \"00.00.0002\", split across two runs each). The fix\n# replaces that leftover placeholder with real sample dates.\n#\n# Walk the document's Fields collection, inspect each field's cached\n# Result text, and for the two placeholders overwrite the underlying Range\n# with the real date. Going through $d.Range(start, end).Text = \"...\" (not\n# the Result range object itself) collapses the old multi-run result into a\n# single run with the new text, which is what the fixed template does.\n\n$d = $word.ActiveDocument\n\n$map = @{\n    \"00.00.0001\" = \"11.11.1970\"\n    \"00.00.0002\" = \"22.22.1970\"\n}\n\n$count = $d.Fields.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $field = $d.Fields.Item($i)\n    $resultRange = $field.Result\n    $currentText = $resultRange.Text\n    if ($map.ContainsKey($currentText)) {\n        $newText = $map[$currentText]\n        $d.Range($resultRange.Start, $resultRange.End).Text = $newText\n    }\n}\n"}
